$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current ("before") values for the columns that move between rows,
# so that source data is not lost before being written to its destination row.
$data = @{}
$data[2] = @{ D = 44637; J = 170; K = 2800; L = 3000; M = 2906; O = "Región Metropolitana"; P = 484 }
$data[3] = @{ D = 44643; J = 90; K = 2800; L = 3000; M = 2911; O = "Región Metropolitana"; P = 485 }
$data[4] = @{ D = 44658; J = 180; K = 2500; L = 3000; M = 2778; O = "Región Metropolitana"; P = 463 }
$data[5] = @{ D = 44659; J = 90; K = 2500; L = 3000; M = 2722; O = "Región Metropolitana"; P = 454 }
$data[6] = @{ D = 44631; J = 110; K = 3000; L = 3500; M = 3273; O = "Provincia de Chacabuco"; P = 546 }
$data[7] = @{ D = 44644; J = 140; K = 2500; L = 3000; M = 2786; O = "Provincia de Chacabuco"; P = 464 }
$data[8] = @{ D = 44672; J = 140; K = 3000; L = 3500; M = 3286; O = "Región Metropolitana"; P = 548 }
$data[9] = @{ D = 44671; J = 150; K = 3500; L = 4000; M = 3733; O = "Región Metropolitana"; P = 622 }
$data[10] = @{ D = 44685; J = 150; K = 3000; L = 3500; M = 3267; O = "Región Metropolitana"; P = 544 }
$data[11] = @{ D = 44630; J = 90; K = 2500; L = 3000; M = 2722; O = "Región Metropolitana"; P = 454 }
$data[12] = @{ D = 44650; J = 130; K = 3000; L = 3500; M = 3308; O = "Región Metropolitana"; P = 551 }

# Apply the permutation: row $dest receives the values previously held by row $src
$mapping = @{}
$mapping[2] = 12
$mapping[3] = 2
$mapping[4] = 5
$mapping[5] = 10
$mapping[6] = 7
$mapping[7] = 11
$mapping[8] = 4
$mapping[9] = 3
$mapping[10] = 6
$mapping[11] = 9
$mapping[12] = 8

foreach ($dest in 2..12) {
    $src = $mapping[$dest]
    $row = $data[$src]
    $ws.Range("D$dest").Value = $row.D
    $ws.Range("J$dest").Value = $row.J
    $ws.Range("K$dest").Value = $row.K
    $ws.Range("L$dest").Value = $row.L
    $ws.Range("M$dest").Value = $row.M
    $ws.Range("O$dest").Value = $row.O
    $ws.Range("P$dest").Value = $row.P
}
